$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A73").Value = "'2024-01-19"
$ws.Range("B73").Value = "11:10:04"
$ws.Range("C73").Value = "Friday"
$ws.Range("D73").Value = "'02"
$ws.Range("E73").Value = 137614
$ws.Range("F73").Value = 140448
$ws.Range("G73").Value = 171050
$ws.Range("H73").Value = 148886
$ws.Range("I73").Value = -1
$ws.Range("J73").Value = 121478
$ws.Range("K73").Value = 223423
$ws.Range("L73").Value = 254823
$ws.Range("M73").Value = 185282
$ws.Range("N73").Value = 110403
$ws.Range("O73").Value = 41357
$ws.Range("P73").Value = 30906
$ws.Range("Q73").Value = 73528
$ws.Range("R73").Value = -1
$ws.Range("S73").Value = 42401
$ws.Range("T73").Value = -1
